# Generate Report for handback
# Adds a new handed-back file entry (785ad82d-54e3-4cf4-8b60-9e02a3834516)
# as an extra row (row 4) on each of the three worksheets:
#   - Overview (summary row)
#   - zh-cn    (detail row for the zh-cn handback)
#   - de-de    (detail row for the de-de handback)

$wb = $excel.ActiveWorkbook

$guid = "785ad82d-54e3-4cf4-8b60-9e02a3834516"
$hash = "e5bc76768abab69fdcac9c9f8767e55711260594"

$mdName  = "$guid.md"
$zhXlfName = "$guid.$hash.zh-cn.xlf"
$deXlfName = "$guid.$hash.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"
$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Overview sheet - one summary row per handed-back file
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8f6c1a2b9d7e4f03c5a6b7d8e9f0a1b2c3d4e5f6/e2e/$mdName",
    [Type]::Missing,
    [Type]::Missing,
    $mdName)
$wsOverview.Range("B4").Value() = $statusInSync
$wsOverview.Range("C4").Value() = $statusInSync

# ---------------------------------------------------------------------------
# zh-cn sheet - detail row for the zh-cn handback
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8f6c1a2b9d7e4f03c5a6b7d8e9f0a1b2c3d4e5f6/e2e/$mdName",
    [Type]::Missing,
    [Type]::Missing,
    $mdName)
$wsZhCn.Range("B4").Value() = $statusInSync
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a2b1c4d5e6f708192a3b4c5d6e7f8091a2b3c4d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName",
    [Type]::Missing,
    [Type]::Missing,
    $zhXlfName)
$wsZhCn.Range("D4").NumberFormat() = $dateFormat
$wsZhCn.Range("D4").Value() = "2016-01-21 02:18:21"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4b3c2d1e0f9a8b7c6d5e4f3a2b1c0d9e8f7a6b5c/e2e/$mdName",
    [Type]::Missing,
    [Type]::Missing,
    $mdName)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5c4d3e2f1a0b9c8d7e6f5a4b3c2d1e0f9a8b7c6d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlfName",
    [Type]::Missing,
    [Type]::Missing,
    $zhXlfName)
$wsZhCn.Range("G4").NumberFormat() = $dateFormat
$wsZhCn.Range("G4").Value() = "2016-01-21 02:19:06"
$wsZhCn.Range("H4").Value() = "Include"

# ---------------------------------------------------------------------------
# de-de sheet - detail row for the de-de handback
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/8f6c1a2b9d7e4f03c5a6b7d8e9f0a1b2c3d4e5f6/e2e/$mdName",
    [Type]::Missing,
    [Type]::Missing,
    $mdName)
$wsDeDe.Range("B4").Value() = $statusInSync
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6d5e4f3a2b1c0d9e8f7a6b5c4d3e2f1a0b9c8d7e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName",
    [Type]::Missing,
    [Type]::Missing,
    $deXlfName)
$wsDeDe.Range("D4").NumberFormat() = $dateFormat
$wsDeDe.Range("D4").Value() = "2016-01-21 02:18:32"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/7e6f5a4b3c2d1e0f9a8b7c6d5e4f3a2b1c0d9e8f/e2e/$mdName",
    [Type]::Missing,
    [Type]::Missing,
    $mdName)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8f7a6b5c4d3e2f1a0b9c8d7e6f5a4b3c2d1e0f9a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlfName",
    [Type]::Missing,
    [Type]::Missing,
    $deXlfName)
$wsDeDe.Range("G4").NumberFormat() = $dateFormat
$wsDeDe.Range("G4").Value() = "2016-01-21 02:19:26"
$wsDeDe.Range("H4").Value() = "Include"

$wb.Save()
